$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force the Price / Volume(1h) data range to Text format so that
# values like "1.00", "5.00" or "68.727.54" are written verbatim instead of
# being auto-coerced into numbers (which would silently drop formatting such
# as trailing zeros or thousands separators used as decimal group dots here).
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '68.727.54'
$ws.Range("E2").Value = '  +2.41%  '

# Row 3
$ws.Range("D3").Value = '2.534.89'
$ws.Range("E3").Value = '  +2.69%  '

# Row 4
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").Value = '594.40'
$ws.Range("E5").Value = '  +2.13%  '

# Row 6
$ws.Range("D6").Value = '178.22'
$ws.Range("E6").Value = '  +2.47%  '

# Row 7
$ws.Range("E7").Value = '  -0.04%  '

# Row 8
$ws.Range("E8").Value = '  +1.20%  '

# Row 9
$ws.Range("D9").Value = '2.534.35'
$ws.Range("E9").Value = '  +2.69%  '

# Row 10
$ws.Range("D10").Value = '0.146'
$ws.Range("E10").Value = '  +6.03%  '

# Row 11
$ws.Range("E11").Value = '  -1.02%  '

# Row 12
$ws.Range("D12").Value = '5.00'
$ws.Range("E12").Value = '  +1.26%  '

# Row 13
$ws.Range("D13").Value = '0.339'
$ws.Range("E13").Value = '  +1.91%  '

# Row 14
$ws.Range("D14").Value = '3.004.56'
$ws.Range("E14").Value = '  +2.80%  '

# Row 15
$ws.Range("D15").Value = '26.11'
$ws.Range("E15").Value = '  +2.98%  '

# Row 16
$ws.Range("D16").Value = '68.471.74'
$ws.Range("E16").Value = '  +2.28%  '

# Row 17
$ws.Range("E17").Value = '  +1.30%  '

# Row 18
$ws.Range("D18").Value = '2.509.23'
$ws.Range("E18").Value = '  +1.54%  '

# Row 19
$ws.Range("D19").Value = '11.12'
$ws.Range("E19").Value = '  +2.32%  '

# Row 20
$ws.Range("E20").Value = '  +0.81%  '

# Row 21
$ws.Range("D21").Value = '353.56'
$ws.Range("E21").Value = '  +1.57%  '

# Row 22
$ws.Range("E22").Value = '  +4.94%  '

# Row 23
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.23%  '

# Row 24
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '71.26'
$ws.Range("E24").Value = '  +2.77%  '

# Row 25
$ws.Range("E25").Value = '  +0.87%  '

# Row 26
$ws.Range("D26").Value = '1.72'
$ws.Range("E26").Value = '  -3.86%  '

# Row 27
$ws.Range("D27").Value = '9.08'
$ws.Range("E27").Value = '  -1.01%  '

# Row 28
$ws.Range("D28").Value = '2.639.11'
$ws.Range("E28").Value = '  +1.70%  '

# Row 29
$ws.Range("E29").Value = '  -0.18%  '

# Row 30
$ws.Range("D30").Value = '515.60'
$ws.Range("E30").Value = '  +3.34%  '

# Row 31
$ws.Range("D31").Value = '0.0₃0901'
$ws.Range("E31").Value = '  +0.21%  '

# Row 32
$ws.Range("D32").Value = '7.82'
$ws.Range("E32").Value = '  +1.45%  '

# Row 33
$ws.Range("D33").Value = '1.26'
$ws.Range("E33").Value = '  +2.61%  '

# Row 34
$ws.Range("E34").Value = '  +1.55%  '

# Row 35
$ws.Range("E35").Value = '  -0.03%  '

# Row 36
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").Value = '164.25'
$ws.Range("E36").Value = '  +1.49%  '

# Row 37
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '0.121'
$ws.Range("E37").Value = '  +0.79%  '

# Row 38
$ws.Range("B38").Value = 'WhiteBITCoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D38").Value = '18.70'
$ws.Range("E38").Value = '  +0.15%  '

# Row 39
$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").Value = '18.44'
$ws.Range("E39").Value = '  +1.92%  '

# Row 40
$ws.Range("D40").Value = '1.32'
$ws.Range("E40").Value = '  +0.31%  '

# Row 41
$ws.Range("D41").Value = '1.77'
$ws.Range("E41").Value = '  +4.96%  '

# Row 42
$ws.Range("E42").Value = '  +0.03%  '

# Row 43
$ws.Range("D43").Value = '4.85'
$ws.Range("E43").Value = '  +0.78%  '

# Row 44
$ws.Range("D44").Value = '0.327'
$ws.Range("E44").Value = '  +0.33%  '

# Row 45
$ws.Range("E45").Value = '  +2.37%  '

# Row 46
$ws.Range("D46").Value = '152.94'
$ws.Range("E46").Value = '  +7.63%  '

# Row 47
$ws.Range("E47").Value = '  +3.01%  '

# Row 48
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0261'
$ws.Range("E48").Value = '  +3.79%  '

# Row 49
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").Value = '0.521'
$ws.Range("E49").Value = '  +2.50%  '

# Row 50
$ws.Range("B50").Value = 'Optimism'
$ws.Range("C50").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("D50").Value = '1.62'
$ws.Range("E50").Value = '  +3.54%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.0742'
$ws.Range("E51").Value = '  +0.41%  '

# Restore the original (unstyled) cell formatting on the data range now that
# the text values are safely stored, so no stray style indices are introduced.
$dataRange.ClearFormats()